$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the row that contains "MGM" in column A (row 302 in the original sheet)
# and delete the entire row, which shifts everything below it up by one,
# matching the commit's removal of the Montgomery, AL entry.
$ws.Rows.Item(302).Delete()
